$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "IrisData" column and its "avg" label
$ws.Range("M1").Value = "IrisData"
$ws.Range("N1").Value = "avg"

# Data values for IrisData column
$ws.Range("M2").Value = 97.33
$ws.Range("M3").Value = 94.67
$ws.Range("M4").Value = 96
$ws.Range("M5").Value = 94.67
$ws.Range("M6").Value = 92
$ws.Range("M7").Value = 96
$ws.Range("M8").Value = 93.33
$ws.Range("M9").Value = 93.33
$ws.Range("M10").Value = 97.33
$ws.Range("M11").Value = 97.33

# Average formula mirroring other columns (e.g. K2 = SUM(J2:J11)/10)
$ws.Range("N2").Formula = "=SUM(M2:M11)/10"

# Update view to reflect new selection/scroll position from the diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("N3").Select()
